# Script 1 - atualização em 2025-09-20 17:07:42Z
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 41 (2023 - Brasil) values
$ws.Cells.Item(41, 3).Value = 35248212
$ws.Cells.Item(41, 4).Value = 64164183
$ws.Cells.Item(41, 5).Value = 4981477

# Update existing row 42 (2023 - Nordeste) values
$ws.Cells.Item(42, 3).Value = 6154165
$ws.Cells.Item(42, 4).Value = 25634101
$ws.Cells.Item(42, 5).Value = 982499

# Row 43 (2023 - Sergipe) remains unchanged

# Add new row 44 (2024 - Brasil)
$ws.Cells.Item(44, 1).Value = 2024
$ws.Cells.Item(44, 2).Value = "Brasil"
$ws.Cells.Item(44, 3).Value = 35743862
$ws.Cells.Item(44, 4).Value = 67313986
$ws.Cells.Item(44, 5).Value = 5409429

# Add new row 45 (2024 - Nordeste)
$ws.Cells.Item(45, 1).Value = 2024
$ws.Cells.Item(45, 2).Value = "Nordeste"
$ws.Cells.Item(45, 3).Value = 6433173
$ws.Cells.Item(45, 4).Value = 26527239
$ws.Cells.Item(45, 5).Value = 1079815

# Add new row 46 (2024 - Sergipe)
$ws.Cells.Item(46, 1).Value = 2024
$ws.Cells.Item(46, 2).Value = "Sergipe"
$ws.Cells.Item(46, 3).Value = 678508
$ws.Cells.Item(46, 4).Value = 192461
$ws.Cells.Item(46, 5).Value = 41268
